# Fill in the newly scraped NFL odds rows (rows 302-327) that were previously
# blank, formatted placeholder rows. Data covers games for dates 1117 and 1120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1117,311,"V","Tennessee",7,7,6,7,27,41.5,40.5,155,20.5),
    @(1117,312,"H","GreenBay",6,0,11,0,17,"pk",3,-175,5),
    @(1120,451,"V","Carolina",0,0,3,0,3,44,41.5,500,20),
    @(1120,452,"H","Baltimore",0,3,0,10,13,12.5,13.5,-700,6),
    @(1120,455,"V","Philadelphia",0,3,0,14,17,10,7,-300,5.5),
    @(1120,456,"H","Indianapolis",7,3,3,3,16,43.5,45.5,250,23.5),
    @(1120,457,"V","Washington",7,13,0,3,23,2.5,3,-160,19),
    @(1120,458,"H","Houston",0,0,3,7,10,40.5,41,140,0.5),
    @(1120,459,"V","NYJets",0,3,0,0,3,40,38.5,155,18.5),
    @(1120,460,"H","NewEngland",0,3,0,7,10,5.5,3.5,-175,2.5),
    @(1120,461,"V","Detroit",3,14,7,7,31,45,44.5,145,23),
    @(1120,462,"H","NYGiants",6,0,0,12,18,4,3,-165,5.5),
    @(1120,463,"V","LARams",7,7,0,6,20,1,39,115,20),
    @(1120,464,"H","NewOrleans",3,7,14,3,27,42.5,1,-135,1),
    @(1120,465,"V","Chicago",7,10,0,7,24,48.5,48,115,1),
    @(1120,466,"H","Atlanta",7,10,7,3,27,4.5,1,-135,23.5),
    @(1120,467,"V","LasVegas",0,7,3,6,22,43.5,41.5,130,21),
    @(1120,468,"H","Denver",7,3,0,6,16,2.5,2.5,-150,0.5),
    @(1120,469,"V","KansasCity",6,7,3,14,30,5.5,5.5,-220,4),
    @(1120,470,"H","LAChargers",10,10,0,7,27,50.5,53,190,27.5),
    @(1120,471,"V","Dallas",10,13,14,3,40,47.5,1,-140,22.5),
    @(1120,472,"H","Minnesota",3,0,0,0,3,"pk",48.5,120,1.5),
    @(1120,473,"V","Cincinnati",10,7,10,10,37,5.5,3.5,-180,3),
    @(1120,474,"H","Pittsburgh",3,17,3,7,30,42.5,39.5,160,21),
    @(1120,477,"N","Cleveland",7,3,0,13,23,48,50.5,290,24.5),
    @(1120,478,"N","Buffalo",3,10,9,9,31,9,9,-350,4)
)

$startRow = 302
$endRow = $startRow + $data.Count - 1

# Rows 305-327 already existed as blank, pre-formatted placeholder rows
# (row-level custom format); strip that row-level formatting now that they
# are becoming real data rows, matching the look of the other data rows
# above them.
$ws.Rows("$startRow`:$endRow").ClearFormats()

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

# The 9 trailing blank placeholder rows that used to sit right after this
# block (old rows 328-336) are no longer part of the data block now that
# rows 302-327 hold real games, so clear any leftover formatting/content in
# that span.
$ws.Range("A328:M336").ClearContents()
$ws.Rows("328:336").ClearFormats()
